$d = $word.ActiveDocument

# Locate the "Author" styled paragraph whose text is exactly "Edison Achalma"
# (this is the author byline right under the title, not the later
# occurrences of the same name in the "Nota de Autores" section).
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Style.NameLocal -eq "Author" -and $p.Range.Text.Trim() -eq "Edison Achalma") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $targetPara = $d.Paragraphs.Item($targetIndex)

    # Insert a new empty paragraph right after the target paragraph by
    # collapsing a range to the end boundary of the paragraph (which sits
    # just past its paragraph mark) and inserting there.
    $endPos = $targetPara.Range.End
    $insertPoint = $d.Range($endPos, $endPos)
    $insertPoint.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Style = "Author"
    $newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
}
